# ClassifierAccuracy.xlsx update: add *_PCA classifier rows, re-sort by Accuracy
# (descending), and leave the selection on B6 as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the eight new rows (top to bottom, so each index below is exactly
#    the live row number at the moment of insertion).
# ---------------------------------------------------------------------------
$ws.Rows(5).Insert()
$ws.Rows(9).Insert()
$ws.Rows(10).Insert()
$ws.Rows(12).Insert()
$ws.Rows(13).Insert()
$ws.Rows(15).Insert()
$ws.Rows(18).Insert()
$ws.Rows(20).Insert()

# ---------------------------------------------------------------------------
# 2. Populate them (Classifier / Accuracy / Hyper Parameters).
# ---------------------------------------------------------------------------
$xlLeft = -4131

$ws.Range("A5").Value = "LinearSVC_PCA"
$ws.Range("B5").Value = 0.76521331458040298
$ws.Range("C5").Value = "{'clf__fit_intercept': True, 'pca__n_components': 100}"
$ws.Range("C5").HorizontalAlignment = $xlLeft

$ws.Range("A9").Value = "MLPClassifier_PCA"
$ws.Range("B9").Value = 0.75527426160337496
$ws.Range("C9").Value = "{'clf__hidden_layer_sizes': (25, 11, 7, 5, 3), 'pca__n_components': 150}"

$ws.Range("A10").Value = "LogisticRegression_PCA"
$ws.Range("B10").Value = 0.75358649789029497
$ws.Range("C10").Value = "{'clf__fit_intercept': True, 'pca__n_components': 100}"

$ws.Range("A12").Value = "RandomForestClassifier_PCA"
$ws.Range("B12").Value = 0.71401781528363795
$ws.Range("C12").Value = "{'clf__max_depth': 10, 'clf__n_estimators': 500, 'pca__n_components': 50}"
$ws.Range("A12").HorizontalAlignment = $xlLeft
$ws.Range("C12").HorizontalAlignment = $xlLeft

$ws.Range("A13").Value = "LinearDiscriminantAnalysis_PCA"
$ws.Range("B13").Value = 0.70961087669948397
$ws.Range("C13").Value = "{'clf__solver': 'lsqr', 'pca__n_components': 100}"
$ws.Range("C13").HorizontalAlignment = $xlLeft

$ws.Range("A15").Value = "KNeighborsClassifier_PCA"
$ws.Range("B15").Value = 0.68469292076887001
$ws.Range("C15").Value = "{'clf__algorithm': 'kd_tree', 'clf__n_neighbors': 15, 'pca__n_components': 50}"

$ws.Range("A18").Value = "RidgeClassifier_PCA"
$ws.Range("B18").Value = 0.66858884200656299
$ws.Range("C18").Value = "{'clf__solver': 'sag', 'pca__n_components': 100}"

$ws.Range("A20").Value = "GaussianNB_PCA"
$ws.Range("B20").Value = 0.57733239568682604
$ws.Range("C20").Value = "{'pca__n_components': 50}"

# ---------------------------------------------------------------------------
# 3. Re-sort the whole table (A2:C22) descending on Accuracy, like the
#    original sorted table, so the sort state / dimension stay accurate.
# ---------------------------------------------------------------------------
$sortRange = $ws.Range("A2:C22")
$keyRange = $ws.Range("B2:B22")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 2)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 4. Restore the cursor position the author left behind.
# ---------------------------------------------------------------------------
$ws.Range("B6").Select()
